$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Constants ---
$xlShiftDown = -4121
$xlPasteFormats = -4122

# 1) Make room for the new worker/period rows. The old row 17 (KEYSI, single
#    period) and the footer rows 22:23 all shift down by 7 rows.
$ws.Range("B17:J23").Insert($xlShiftDown)

# 2) The freshly inserted rows 17:23 come in blank/unformatted - clone the
#    "interior" row look (row 16) into them so they pick up the same
#    borders/fill/number-format as the rest of the table body.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J23").PasteSpecial($xlPasteFormats)

# 3) Populate the new worker rows.
# New worker: ALEXIDES PALACIN DE LA HOZ
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "9297160"
$ws.Range("D17").Value = "ALEXIDES PALACIN DE LA HOZ"
$ws.Range("E17").Value = "2306"
$ws.Range("F17").Value = 59000
$ws.Range("G17").Value = 1450000

# Existing worker KEYSI NINOSKA ALDAVES OCHOA, moved up into the normal body style
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143377852"
$ws.Range("D18").Value = "KEYSI NINOSKA ALDAVES OCHOA"
$ws.Range("E18").Value = "2010"
$ws.Range("F18").Value = 19200
$ws.Range("G18").Value = 1200000

# New worker: CARLOS ANDRES CAMACHO CABARCAS - six overdue periods (2502-2507)
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1001971310"
$ws.Range("D19").Value = "CARLOS ANDRES CAMACHO CABARCAS"
$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 52000
$ws.Range("G19").Value = 1300000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1001971310"
$ws.Range("D20").Value = "CARLOS ANDRES CAMACHO CABARCAS"
$ws.Range("E20").Value = "2506"
$ws.Range("F20").Value = 52000
$ws.Range("G20").Value = 1300000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1001971310"
$ws.Range("D21").Value = "CARLOS ANDRES CAMACHO CABARCAS"
$ws.Range("E21").Value = "2505"
$ws.Range("F21").Value = 52000
$ws.Range("G21").Value = 1300000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1001971310"
$ws.Range("D22").Value = "CARLOS ANDRES CAMACHO CABARCAS"
$ws.Range("E22").Value = "2504"
$ws.Range("F22").Value = 52000
$ws.Range("G22").Value = 1300000

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1001971310"
$ws.Range("D23").Value = "CARLOS ANDRES CAMACHO CABARCAS"
$ws.Range("E23").Value = "2503"
$ws.Range("F23").Value = 52000
$ws.Range("G23").Value = 1300000

# 4) Row 24 still holds the old (now stale) KEYSI record but kept the special
#    "last row" bottom-border style - overwrite it with the final Carlos period.
$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1001971310"
$ws.Range("D24").Value = "CARLOS ANDRES CAMACHO CABARCAS"
$ws.Range("E24").Value = "2502"
$ws.Range("F24").Value = 52000
$ws.Range("G24").Value = 1300000

# 5) Refresh the summary figures.
$ws.Range("E11").Value = 407756
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 9
